$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53: No Accounting for Waste / Enchanted Electrum Ink (Item ID 5479)
$ws.Range("H53").Value = 841.1429000000001
$ws.Range("I53").Value = 1567.7142
$ws.Range("J53").Value = 114.57143
$ws.Range("K53").Value = 1567.7142
$ws.Range("L53").Value = 114.57143
$ws.Range("M53").Value = -930.7141999999999
$ws.Range("N53").Value = -1388.57143

# Row 70: Consecrating Congregation / Holy Water (Item ID 12604)
$ws.Range("H70").Value = 1294.15
$ws.Range("I70").Value = 885
$ws.Range("J70").Value = 1566.9166
$ws.Range("K70").Value = 2655
$ws.Range("L70").Value = 4700.7498
$ws.Range("M70").Value = -2385
$ws.Range("N70").Value = -5240.7498

# Row 73: Curbing the Contagion (L) / Holy Water (Item ID 12604)
$ws.Range("H73").Value = 1294.15
$ws.Range("I73").Value = 885
$ws.Range("J73").Value = 1566.9166
$ws.Range("K73").Value = 2655
$ws.Range("L73").Value = 4700.7498
$ws.Range("M73").Value = -1719
$ws.Range("N73").Value = -6572.7498

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone (Item ID 44013)
$ws.Range("H137").Value = 9091871
$ws.Range("I137").Value = 614.2857
$ws.Range("J137").Value = 13334458
$ws.Range("K137").Value = 1842.8571
$ws.Range("L137").Value = 40003374
$ws.Range("M137").Value = 707.1428999999998
$ws.Range("N137").Value = -40008474

# Row 138: All-night Crafting / Cunning Craftsman's Tisane (Item ID 44169)
$ws.Range("H138").Value = 5466408
$ws.Range("I138").Value = 11112267
$ws.Range("J138").Value = 2674.1936
$ws.Range("K138").Value = 33336801
$ws.Range("L138").Value = 8022.5808
$ws.Range("M138").Value = -33331661
$ws.Range("N138").Value = -18302.5808

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot (Item ID 43999)
$ws.Range("H61").Value = 19233178
$ws.Range("I61").Value = 20002106
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 20002106
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -20001894
$ws.Range("N61").Value = -10424

# Row 74: As the Bolt Flies / Titanium Nugget (Item ID 44000)
$ws.Range("H74").Value = 15155166
$ws.Range("I74").Value = 20835814
$ws.Range("J74").Value = 6769.778
$ws.Range("K74").Value = 20835814
$ws.Range("L74").Value = 6769.778
$ws.Range("M74").Value = -20834940
$ws.Range("N74").Value = -8517.778

# Row 77: Heavy Metal Banned (L) / Titanium Nugget (Item ID 44000)
$ws.Range("H77").Value = 15155166
$ws.Range("I77").Value = 20835814
$ws.Range("J77").Value = 6769.778
$ws.Range("K77").Value = 104179070
$ws.Range("L77").Value = 33848.89
$ws.Range("M77").Value = -104174702
$ws.Range("N77").Value = -42584.89

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot (Item ID 43999)
$ws.Range("H136").Value = 19233178
$ws.Range("I136").Value = 20002106
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 60006318
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -60003768
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot (Item ID 43998)
$ws.Range("H134").Value = 2444.2
$ws.Range("I134").Value = 1373.5778
$ws.Range("J134").Value = 5656.067
$ws.Range("K134").Value = 4120.7334
$ws.Range("L134").Value = 16968.201
$ws.Range("M134").Value = -1585.7334
$ws.Range("N134").Value = -22038.201

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber (Item ID 44023)
$ws.Range("H31").Value = 14450767
$ws.Range("I31").Value = 4006626.8
$ws.Range("K31").Value = 4006626.8
$ws.Range("M31").Value = -4006331.8

# Row 34: Armoires of the Rich and Famous / Walnut Lumber (Item ID 44023)
$ws.Range("H34").Value = 14450767
$ws.Range("I34").Value = 4006626.8
$ws.Range("K34").Value = 4006626.8
$ws.Range("M34").Value = -4006424.8

# Row 62: Splinter in the Sewers / Cedar Lumber (Item ID 12580)
$ws.Range("H62").Value = 2077.64
$ws.Range("I62").Value = 1988.0555
$ws.Range("J62").Value = 2308
$ws.Range("K62").Value = 1988.0555
$ws.Range("L62").Value = 2308
$ws.Range("M62").Value = -1364.0555
$ws.Range("N62").Value = -3556

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber (Item ID 12580)
$ws.Range("H65").Value = 2077.64
$ws.Range("I65").Value = 1988.0555
$ws.Range("J65").Value = 2308
$ws.Range("K65").Value = 9940.2775
$ws.Range("L65").Value = 11540
$ws.Range("M65").Value = -6820.2775
$ws.Range("N65").Value = -17780

# Row 132: Hull Lotta Damage / Ginseng Lumber (Item ID 44019)
$ws.Range("H132").Value = 6411538
$ws.Range("I132").Value = 7693392
$ws.Range("J132").Value = 2268.1538
$ws.Range("K132").Value = 23080176
$ws.Range("L132").Value = 6804.4614
$ws.Range("M132").Value = -23077646
$ws.Range("N132").Value = -11864.4614

# Row 134: Wood You Be Quiet / Ceiba Lumber (Item ID 44020)
$ws.Range("H134").Value = 291524.56
$ws.Range("I134").Value = 1090.2458
$ws.Range("J134").Value = 1135167.1
$ws.Range("K134").Value = 3270.7374
$ws.Range("L134").Value = 3405501.3
$ws.Range("M134").Value = -735.7374
$ws.Range("N134").Value = -3410571.3

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service / Frantoio Oil (Item ID 27838)
$ws.Range("H107").Value = 698.2
$ws.Range("I107").Value = 118.75
$ws.Range("J107").Value = 970.8823
$ws.Range("K107").Value = 356.25
$ws.Range("L107").Value = 2912.6469
$ws.Range("M107").Value = 1563.75
$ws.Range("N107").Value = -6752.6469

# Row 132: More Mezcal / Cooking Mezcal (Item ID 43972)
$ws.Range("H132").Value = 1597.6364
$ws.Range("I132").Value = 551
$ws.Range("J132").Value = 2195.7144
$ws.Range("K132").Value = 4959
$ws.Range("L132").Value = 19761.4296
$ws.Range("M132").Value = -2429
$ws.Range("N132").Value = -24821.4296

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar / Lar Ingot (Item ID 44008)
$ws.Range("H132").Value = 3411.4043
$ws.Range("I132").Value = 2521.6667
$ws.Range("J132").Value = 5508.643
$ws.Range("K132").Value = 7565.000100000001
$ws.Range("L132").Value = 16525.929
$ws.Range("M132").Value = -5035.000100000001
$ws.Range("N132").Value = -21585.929

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather (Item ID 36249)
$ws.Range("H7").Value = 8344.666999999999
$ws.Range("I7").Value = 30002
$ws.Range("J7").Value = 5637.5
$ws.Range("K7").Value = 30002
$ws.Range("L7").Value = 5637.5
$ws.Range("M7").Value = -29890
$ws.Range("N7").Value = -5861.5

# Row 93: Hide to Go Seek / Gagana Leather (Item ID 19993)
$ws.Range("H93").Value = 981.8125
$ws.Range("I93").Value = 953.9
$ws.Range("J93").Value = 1028.3334
$ws.Range("K93").Value = 953.9
$ws.Range("L93").Value = 1028.3334
$ws.Range("M93").Value = 294.1
$ws.Range("N93").Value = -3524.3334

# Row 126: Battered Books / Saiga Leather (Item ID 36249)
$ws.Range("H126").Value = 8344.666999999999
$ws.Range("I126").Value = 30002
$ws.Range("J126").Value = 5637.5
$ws.Range("K126").Value = 90006
$ws.Range("L126").Value = 16912.5
$ws.Range("M126").Value = -87536
$ws.Range("N126").Value = -21852.5

# Row 132: Tenets of Tanning / Silver Lobo Leather (Item ID 44058)
$ws.Range("H132").Value = 8936490
$ws.Range("I132").Value = 5312.6665
$ws.Range("J132").Value = 25012608
$ws.Range("K132").Value = 15937.9995
$ws.Range("L132").Value = 75037824
$ws.Range("M132").Value = -13407.9995
$ws.Range("N132").Value = -75042884

# Row 136: Respect for Br'aax / Br'aax Leather (Item ID 44060)
$ws.Range("H136").Value = 7248808.5
$ws.Range("I136").Value = 8475498
$ws.Range("J136").Value = 11341.5
$ws.Range("K136").Value = 25426494
$ws.Range("L136").Value = 34024.5
$ws.Range("M136").Value = -25423944
$ws.Range("N136").Value = -39124.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth (Item ID 44029)
$ws.Range("H132").Value = 2265.5833
$ws.Range("I132").Value = 1950.4482
$ws.Range("J132").Value = 3571.1428
$ws.Range("K132").Value = 5851.3446
$ws.Range("L132").Value = 10713.4284
$ws.Range("M132").Value = -3321.3446
$ws.Range("N132").Value = -15773.4284

# Row 136: Weaving the Envelope / Sarcenet Cloth (Item ID 44031)
$ws.Range("H136").Value = 835.431
$ws.Range("I136").Value = 675.51166
$ws.Range("J136").Value = 1293.8667
$ws.Range("K136").Value = 2026.53498
$ws.Range("L136").Value = 3881.6001
$ws.Range("M136").Value = 523.4650200000001
$ws.Range("N136").Value = -8981.6001
